$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Update existing timer-divider values (Fclk stays 1,000,000 Hz)
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 64          # Timer1 divider: 1024 -> 64
$ws.Range("B7").Value = 200         # Cycle duration (mS): 500 -> 200
$ws.Range("B10").Value = 8.4        # Packet duration (ms): 36 -> 8.4

# Row 11 used to hold a literal "tics per packet" value; it becomes a
# calculated ("Check Cell"-bordered => "Calculation") row driven by a
# formula, matching row 8's pattern (B10/B6).
$ws.Range("A11").Style = "Calculation"
$ws.Range("B11").Formula = "=B10/B6"
$ws.Range("B11").Style = "Calculation"

# ---------------------------------------------------------------------
# 2. Remove the old "Address / WhenTransmit" block (rows 11 D-col, 12-14)
# ---------------------------------------------------------------------
$ws.Range("D11").ClearContents()
$ws.Range("A12:D14").Clear()

# ---------------------------------------------------------------------
# 3. New "Address" block (rows 16-19) and "Time" block (rows 21-22).
#    Shared-string table order matters for a byte-faithful rebuild, so
#    new labels are introduced in the same order the original author
#    typed them in: Address, WhenTransmit, Time, Multiplier, Addr N...
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Address"
$ws.Range("A16").Style = "Good"
$ws.Range("B16").Value = 1
$ws.Range("B16").Style = "Good"

$ws.Range("A18").Value = "WhenTransmit"
$ws.Range("A18").Style = "Calculation"
$ws.Range("C18").Value = "tics"

$ws.Range("A19").Value = "WhenTransmit"
$ws.Range("C19").Value = "ms"

$ws.Range("A21").Value = "Time"
$ws.Range("A21").Style = "Good"
$ws.Range("B21").Value = 7.4
$ws.Range("B21").Style = "Good"
$ws.Range("C21").Value = "ms"

$ws.Range("A22").Value = "Time"
$ws.Range("A22").Style = "Calculation"
$ws.Range("B22").Formula = "=B21/B6"
$ws.Range("B22").Style = "Calculation"
$ws.Range("C22").Value = "tics"

$ws.Range("A17").Value = "Multiplier"
$ws.Range("A17").Style = "Good"
$ws.Range("B17").Value = 128
$ws.Range("B17").Style = "Good"

# Named ranges referenced by formulas below - defined now that B16/B17 exist.
$wb.Names.Add("Addr", $ws.Range("B16"))
$wb.Names.Add("Multi", $ws.Range("B17"))

$ws.Range("B18").Formula = "=Addr*B17+B11"
$ws.Range("B18").Style = "Calculation"

$ws.Range("B19").Formula = "=B18*B6"

# ---------------------------------------------------------------------
# 5. "Addr N" multiplication table (rows 25-34)
# ---------------------------------------------------------------------
$addrLabels = @("Addr 1","Addr 2","Addr 3","Addr 4","Addr 5","Addr 6","Addr 7","Addr 8","Addr 9","Addr 10")
for ($i = 0; $i -lt 10; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value = $addrLabels[$i]
    $ws.Cells.Item($r, 2).Value = $i + 1
    $ws.Cells.Item($r, 3).Formula = "=B$r*Multi"
}

# ---------------------------------------------------------------------
# 6. Sheet view: selection + (best effort) scrolled-into-view top row
# ---------------------------------------------------------------------
$ws.Range("B18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# 7. Page setup (paper size + orientation)
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 8. Workbook window geometry (best effort - cosmetic only)
# ---------------------------------------------------------------------
try {
    $excel.ActiveWindow.Width = 19080
    $excel.ActiveWindow.Height = 8475
} catch {
}
